# Auto-generated Excel COM-interop script to apply cryptos list update
# (commit: "Updated cryptos list on Sun Jun 30 05:11:11 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.757.54"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "3.361.22"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.67"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E11").Value = "  -4.20%  "
$ws.Range("D12").Value = "3.936.84"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "3.365.77"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "60.877.60"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.08"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.47"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.87"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.58"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +7.14%  "
$ws.Range("E26").Value = "  -7.06%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  -4.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.81"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -6.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.92"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "167.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "3.395.97"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.43"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.16"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -10.52%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.768"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.10"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.448.03"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -3.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0256"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("E51").Value = "  -3.24%  "
